$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (sheet ALC)
$ws.Range("H2").Value = 473.2857
$ws.Range("I2").Value = 604.8182
$ws.Range("J2").Value = 388.17648
$ws.Range("K2").Value = 604.8182
$ws.Range("L2").Value = 388.17648
$ws.Range("M2").Value = -491.8182
$ws.Range("N2").Value = -614.1764800000001

# Row 4 (sheet ALC)
$ws.Range("H4").Value = 250.90909
$ws.Range("I4").Value = 177.75
$ws.Range("J4").Value = 446
$ws.Range("K4").Value = 177.75
$ws.Range("L4").Value = 446
$ws.Range("M4").Value = -63.75

# Row 33 (sheet ALC)
$ws.Range("H33").Value = 783.48
$ws.Range("I33").Value = 463
$ws.Range("J33").Value = 2466
$ws.Range("K33").Value = 463
$ws.Range("L33").Value = 2466
$ws.Range("M33").Value = -234

# Row 40 (sheet ALC)
$ws.Range("H40").Value = 2574.25
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 2599
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2599
$ws.Range("M40").Value = -2325

# Row 41 (sheet ALC)
$ws.Range("H41").Value = 542.86664
$ws.Range("I41").Value = 441.76923
$ws.Range("J41").Value = 1200
$ws.Range("K41").Value = 441.76923
$ws.Range("L41").Value = 1200
$ws.Range("M41").Value = -1.769229999999993

# Row 53 (sheet ALC)
$ws.Range("H53").Value = 562.92
$ws.Range("I53").Value = 495.93332
$ws.Range("J53").Value = 663.4
$ws.Range("K53").Value = 495.93332
$ws.Range("L53").Value = 663.4
$ws.Range("M53").Value = 141.06668
$ws.Range("N53").Value = -1937.4

# Row 62 (sheet ALC)
$ws.Range("H62").Value = 3676.2
$ws.Range("I62").Value = 3373.25
$ws.Range("J62").Value = 4888
$ws.Range("K62").Value = 3373.25
$ws.Range("L62").Value = 4888
$ws.Range("M62").Value = -2749.25

# Row 65 (sheet ALC)
$ws.Range("H65").Value = 3676.2
$ws.Range("I65").Value = 3373.25
$ws.Range("J65").Value = 4888
$ws.Range("K65").Value = 16866.25
$ws.Range("L65").Value = 24440
$ws.Range("M65").Value = -13746.25

# Row 70 (sheet ALC)
$ws.Range("H70").Value = 1397
$ws.Range("I70").Value = 1308.1666
$ws.Range("J70").Value = 1441.4166
$ws.Range("K70").Value = 3924.4998
$ws.Range("L70").Value = 4324.2498
$ws.Range("M70").Value = -3654.4998

# Row 73 (sheet ALC)
$ws.Range("H73").Value = 1397
$ws.Range("I73").Value = 1308.1666
$ws.Range("J73").Value = 1441.4166
$ws.Range("K73").Value = 3924.4998
$ws.Range("L73").Value = 4324.2498
$ws.Range("M73").Value = -2988.4998

# Row 113 (sheet ALC)
$ws.Range("H113").Value = 27500
$ws.Range("I113").Value = 27500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 27500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -24246
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 4 (sheet ARM)
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()

# Row 22 (sheet ARM)
$ws.Range("H22").Value = 19000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 19000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 19000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -19598

# Row 45 (sheet ARM)
$ws.Range("H45").Value = 3319.9092
$ws.Range("I45").Value = 1129.75
$ws.Range("J45").Value = 4571.4287
$ws.Range("K45").Value = 1129.75
$ws.Range("L45").Value = 4571.4287
$ws.Range("M45").Value = -752.75

# Row 50 (sheet ARM)
$ws.Range("H50").Value = 6078.727
$ws.Range("I50").Value = 7052.2856
$ws.Range("J50").Value = 4375
$ws.Range("K50").Value = 7052.2856
$ws.Range("L50").Value = 4375
$ws.Range("M50").Value = -6338.2856
$ws.Range("N50").Value = -5803

# Row 97 (sheet ARM)
$ws.Range("H97").Value = 1303.625
$ws.Range("I97").Value = 737.25
$ws.Range("J97").Value = 3002.75
$ws.Range("K97").Value = 737.25
$ws.Range("L97").Value = 3002.75
$ws.Range("M97").Value = -241.25
$ws.Range("N97").Value = -3994.75

# Row 122 (sheet ARM)
$ws.Range("H122").Value = 2575.889
$ws.Range("I122").Value = 2176.077
$ws.Range("J122").Value = 3615.4
$ws.Range("K122").Value = 6528.231000000001
$ws.Range("L122").Value = 10846.2
$ws.Range("M122").Value = -4078.231000000001

# Row 134 (sheet ARM)
$ws.Range("H134").Value = 80000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 80000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (sheet CRP)
$ws.Range("H7").Value = 1390.6072
$ws.Range("I7").Value = 681.2632
$ws.Range("J7").Value = 2888.111
$ws.Range("K7").Value = 681.2632
$ws.Range("L7").Value = 2888.111
$ws.Range("M7").Value = -568.2632

# Row 22 (sheet CRP)
$ws.Range("H22").Value = 689.6
$ws.Range("I22").Value = 699.5
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 699.5
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -349.5
$ws.Range("N22").Value = -1350

# Row 29 (sheet CRP)
$ws.Range("H29").Value = 6500
$ws.Range("I29").Value = 6500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 6500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -6207
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 19 (sheet CUL)
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 41 (sheet CUL)
$ws.Range("H41").Value = 250.5
$ws.Range("I41").Value = 250.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 751.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -413.5

# Row 131 (sheet CUL)
$ws.Range("H131").Value = 1067.7142
$ws.Range("I131").Value = 799.6667
$ws.Range("J131").Value = 1099.88
$ws.Range("K131").Value = 2399.0001
$ws.Range("L131").Value = 3299.64
$ws.Range("M131").Value = 2640.9999
$ws.Range("N131").Value = -13379.64

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (sheet GSM)
$ws.Range("H2").Value = 639.26666
$ws.Range("I2").Value = 507.1111
$ws.Range("J2").Value = 837.5
$ws.Range("K2").Value = 507.1111
$ws.Range("L2").Value = 837.5
$ws.Range("M2").Value = -394.1111

# Row 80 (sheet GSM)
$ws.Range("H80").Value = 2626.5
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3253
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 3253
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -5249

# Row 83 (sheet GSM)
$ws.Range("H83").Value = 2626.5
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3253
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 16265
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -26249

# Row 122 (sheet GSM)
$ws.Range("H122").Value = 1322.3334
$ws.Range("I122").Value = 1364.125
$ws.Range("J122").Value = 988
$ws.Range("K122").Value = 4092.375
$ws.Range("L122").Value = 2964
$ws.Range("M122").Value = -1642.375

# Row 132 (sheet GSM)
$ws.Range("H132").Value = 7424.125
$ws.Range("I132").Value = 7400.6665
$ws.Range("J132").Value = 7494.5
$ws.Range("K132").Value = 22201.9995
$ws.Range("L132").Value = 22483.5
$ws.Range("M132").Value = -19671.9995
$ws.Range("N132").Value = -27543.5

# Row 134 (sheet GSM)
$ws.Range("H134").Value = 25663
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 25663
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 76989
$ws.Range("N134").Value = -82059

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (sheet LTW)
$ws.Range("H22").Value = 2596.75
$ws.Range("I22").Value = 1909.8
$ws.Range("J22").Value = 3741.6667
$ws.Range("K22").Value = 1909.8
$ws.Range("L22").Value = 3741.6667
$ws.Range("M22").Value = -1614.8
$ws.Range("N22").Value = -4331.6667

# Row 27 (sheet LTW)
$ws.Range("H27").Value = 2596.75
$ws.Range("I27").Value = 1909.8
$ws.Range("J27").Value = 3741.6667
$ws.Range("K27").Value = 1909.8
$ws.Range("L27").Value = 3741.6667
$ws.Range("M27").Value = -1802.8
$ws.Range("N27").Value = -3955.6667

# Row 46 (sheet LTW)
$ws.Range("H46").Value = 1472.9474
$ws.Range("I46").Value = 1271.091
$ws.Range("J46").Value = 1750.5
$ws.Range("K46").Value = 1271.091
$ws.Range("L46").Value = 1750.5
$ws.Range("M46").Value = -1083.091

# Row 68 (sheet LTW)
$ws.Range("H68").Value = 1895
$ws.Range("I68").Value = 1895
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1895
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1146
$ws.Range("N68").ClearContents()

# Row 71 (sheet LTW)
$ws.Range("H71").Value = 1895
$ws.Range("I71").Value = 1895
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9475
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5731
$ws.Range("N71").ClearContents()

# Row 98 (sheet LTW)
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (sheet WVR)
$ws.Range("H100").Value = 5809759.5
$ws.Range("I100").Value = 8712965
$ws.Range("J100").Value = 3349.25
$ws.Range("K100").Value = 17425930
$ws.Range("L100").Value = 6698.5
$ws.Range("M100").Value = -17425389
$ws.Range("N100").Value = -7780.5

# Row 113 (sheet WVR)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# Row 132 (sheet WVR)
$ws.Range("H132").Value = 2871.3076
$ws.Range("I132").Value = 2865.7273
$ws.Range("J132").Value = 2902
$ws.Range("K132").Value = 8597.1819
$ws.Range("L132").Value = 8706
$ws.Range("M132").Value = -6067.1819
$ws.Range("N132").Value = -13766
